$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the updated row 4 (POL/POD uppercased, surcharges added) plus the
# newly added rows 5-9.
# Columns: A=POL, B=POD, C=Container, D=Freight USD, E=OTHC AUD,
#          F=DOC AUD, G=CMR AUD, H=AMS USD, I=LSS USD, J=DTHC, K=Free Time
$rows = @(
    @{ Row = 4; POL = "MELBOURNE"; POD = "TAICHUNG"; Freight = "500" },
    @{ Row = 5; POL = "SYDNEY";    POD = "TAICHUNG"; Freight = "500" },
    @{ Row = 6; POL = "BRISBANE";  POD = "TAICHUNG"; Freight = "500" },
    @{ Row = 7; POL = "AU";        POD = "TAICHUNG"; Freight = "500" },
    @{ Row = 8; POL = "MELBOURNE"; POD = "SHANGHAI"; Freight = "500" },
    @{ Row = 9; POL = "MELBOURNE"; POD = "TOKYO";    Freight = 800 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.POL
    $ws.Cells.Item($rowNum, 2).Value = $r.POD
    $ws.Cells.Item($rowNum, 3).Value = "40GP"

    # Freight USD (column D) - text for most rows, numeric for row 9
    if ($r.Freight -is [string]) {
        $ws.Cells.Item($rowNum, 4).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 4).Value = $r.Freight
        $ws.Cells.Item($rowNum, 4).Style = "Normal"
    } else {
        $ws.Cells.Item($rowNum, 4).Value = $r.Freight
    }

    # OTHC AUD (column E) stays textual "400"
    $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 5).Value = "400"
    $ws.Cells.Item($rowNum, 5).Style = "Normal"

    # DOC AUD, CMR AUD, AMS USD, LSS USD - numeric values
    $ws.Cells.Item($rowNum, 6).Value = 120
    $ws.Cells.Item($rowNum, 7).Value = 20
    $ws.Cells.Item($rowNum, 8).Value = 30
    $ws.Cells.Item($rowNum, 9).Value = 70

    $ws.Cells.Item($rowNum, 10).Value = "Collect"
    $ws.Cells.Item($rowNum, 11).Value = "14 Days"
}
